# Applies the "Add files via upload" re-upload of game_data_input.xlsx:
# the header row (A1:AF1) and the label columns A/C/D/E are unchanged;
# every KPI measure in rows 2-4, plus the two category labels (M, AC)
# and a couple of state codes (B3, B4), were replaced with a refreshed
# data pull. Writing the new literal text values first (in the same
# order they were newly introduced) lets Excel's own shared-string
# table get rebuilt dropping now-unused strings ("Slurpee", "Bakery",
# "Redbull", "Carbonated Soft Drinks", "Chicken") and appending the new
# ones ("CA", "Grill", "CSD", "Energy Drinks", "Cigarettes", "MI") in
# the same order the refreshed workbook has them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New text values, in first-use order ---
$ws.Range("B3").Value = "CA"
$ws.Range("M2").Value = "Grill"
$ws.Range("AC2").Value = "CSD"
$ws.Range("M3").Value = "Energy Drinks"
$ws.Range("AC3").Value = "Chips"
$ws.Range("M4").Value = "Cigarettes"
$ws.Range("AC4").Value = "CSD"
$ws.Range("B4").Value = "MI"

# --- Row 2 (store 10017) ---
$ws.Range("F2").Value = 1735
$ws.Range("G2").Value = 10577
$ws.Range("H2").Value = 3697
$ws.Range("I2").Value = 0.349532003403611
$ws.Range("J2").Value = 291
$ws.Range("K2").Value = 2780
$ws.Range("L2").Value = 742
$ws.Range("N2").Value = 348
$ws.Range("O2").Value = 0.20057636887608
$ws.Range("P2").Value = 42.1436781609195
$ws.Range("Q2").Value = 870
$ws.Range("R2").Value = 0.501440922190201
$ws.Range("S2").Value = 5615
$ws.Range("T2").Value = 6.45402298850574
$ws.Range("U2").Value = 155
$ws.Range("V2").Value = 0.0893371757925072
$ws.Range("W2").Value = 22509.9199436716
$ws.Range("X2").Value = 145.225289959172
$ws.Range("Y2").Value = 219
$ws.Range("Z2").Value = 0.126224783861671
$ws.Range("AA2").Value = 2919
$ws.Range("AB2").Value = 13.3287671232876
$ws.Range("AD2").Value = 938
$ws.Range("AE2").Value = 545
$ws.Range("AF2").Value = 0.581023454157782

# --- Row 3 (store 10022) ---
$ws.Range("F3").Value = 1654
$ws.Range("G3").Value = 14839
$ws.Range("H3").Value = 4466
$ws.Range("I3").Value = 0.300963676797627
$ws.Range("J3").Value = 633
$ws.Range("K3").Value = 2376
$ws.Range("L3").Value = 1100
$ws.Range("N3").Value = 499
$ws.Range("O3").Value = 0.301692865779927
$ws.Range("P3").Value = 60.5330661322645
$ws.Range("Q3").Value = 933
$ws.Range("R3").Value = 0.564087061668681
$ws.Range("S3").Value = 11268
$ws.Range("T3").Value = 12.0771704180064
$ws.Range("U3").Value = 180
$ws.Range("V3").Value = 0.108827085852478
$ws.Range("W3").Value = 41429.3899835161
$ws.Range("X3").Value = 230.163277686201
$ws.Range("Y3").Value = 436
$ws.Range("Z3").Value = 0.263603385731559
$ws.Range("AA3").Value = 9138
$ws.Range("AB3").Value = 20.9587155963302
$ws.Range("AD3").Value = 1380
$ws.Range("AE3").Value = 1085
$ws.Range("AF3").Value = 0.786231884057971

# --- Row 4 (store 10056) ---
$ws.Range("F4").Value = 2168
$ws.Range("G4").Value = 17917
$ws.Range("H4").Value = 7165
$ws.Range("I4").Value = 0.399899536752804
$ws.Range("J4").Value = 215
$ws.Range("K4").Value = 1733
$ws.Range("L4").Value = 1193
$ws.Range("N4").Value = 839
$ws.Range("O4").Value = 0.386992619926199
$ws.Range("P4").Value = 60.8641239570917
$ws.Range("Q4").Value = 1479
$ws.Range("R4").Value = 0.682195571955719
$ws.Range("S4").Value = 20751
$ws.Range("T4").Value = 14.0304259634888
$ws.Range("U4").Value = 239
$ws.Range("V4").Value = 0.110239852398523
$ws.Range("W4").Value = 57506.6100597493
$ws.Range("X4").Value = 240.61343121234
$ws.Range("Y4").Value = 692
$ws.Range("Z4").Value = 0.319188191881918
$ws.Range("AA4").Value = 15371
$ws.Range("AB4").Value = 22.2124277456647
$ws.Range("AD4").Value = 1375
$ws.Range("AE4").Value = 782
$ws.Range("AF4").Value = 0.568727272727272

# --- Window/view state: the saved file now has the cursor on AF2 with
# the grid scrolled so column AA is the left-most visible column.
$excel.ActiveWindow.ScrollColumn = 27
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AF2").Select()
